$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 4384
$ws.Range("E2").Value = 67
$ws.Range("F2").Value = 67
$ws.Range("G2").Value = 50
$ws.Range("H2").Value = 26
$ws.Range("I2").Value = 26
$ws.Range("K2").Value = 2283
$ws.Range("L2").Value = 587
$ws.Range("M2").Value = 1696
$ws.Range("N2").Value = 1696
$ws.Range("P2").Value = 64
$ws.Range("Q2").Value = -110
$ws.Range("R2").Value = -81
$ws.Range("S2").Value = -35
$ws.Range("T2").Value = 61
$ws.Range("U2").Value = -171
$ws.Range("W2").Value = 1.54
$ws.Range("X2").Value = 0.59
$ws.Range("Y2").Value = 1.52
$ws.Range("Z2").Value = 1.1
$ws.Range("AA2").Value = 34.63
$ws.Range("AB2").Value = 2531.29
$ws.Range("AC2").Value = 146
$ws.Range("AD2").Value = 127.5
$ws.Range("AE2").Value = 9780
$ws.Range("AF2").Value = 1.9
$ws.Range("AG2").Value = 7
$ws.Range("AH2").Value = 0.04
$ws.Range("AI2").Value = 4.87
$ws.Range("AJ2").Value = 17588764
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("V2").ClearContents()

# Row 3
$ws.Range("D3").Value = 4079
$ws.Range("E3").Value = 177
$ws.Range("F3").Value = 177
$ws.Range("G3").Value = 204
$ws.Range("H3").Value = 156
$ws.Range("I3").Value = 156
$ws.Range("K3").Value = 2486
$ws.Range("L3").Value = 632
$ws.Range("M3").Value = 1854
$ws.Range("N3").Value = 1854
$ws.Range("P3").Value = 70
$ws.Range("Q3").Value = 336
$ws.Range("R3").Value = -148
$ws.Range("S3").Value = -3
$ws.Range("T3").Value = 29
$ws.Range("U3").Value = 307
$ws.Range("W3").Value = 4.34
$ws.Range("X3").Value = 3.82
$ws.Range("Y3").Value = 8.779999999999999
$ws.Range("Z3").Value = 6.54
$ws.Range("AA3").Value = 34.07
$ws.Range("AB3").Value = 2512.43
$ws.Range("AC3").Value = 886
$ws.Range("AD3").Value = 23.13
$ws.Range("AE3").Value = 10696
$ws.Range("AF3").Value = 1.92
$ws.Range("AG3").Value = 239
$ws.Range("AH3").Value = 1.17
$ws.Range("AI3").Value = 26.46
$ws.Range("AJ3").Value = 17588764
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("V3").ClearContents()

# Row 4
$ws.Range("D4").Value = 4346
$ws.Range("E4").Value = 243
$ws.Range("F4").Value = 243
$ws.Range("G4").Value = 269
$ws.Range("H4").Value = 180
$ws.Range("I4").Value = 180
$ws.Range("K4").Value = 2779
$ws.Range("L4").Value = 797
$ws.Range("M4").Value = 1981
$ws.Range("N4").Value = 1981
$ws.Range("P4").Value = 77
$ws.Range("Q4").Value = 272
$ws.Range("R4").Value = -44
$ws.Range("S4").Value = -43
$ws.Range("T4").Value = 30
$ws.Range("U4").Value = 242
$ws.Range("W4").Value = 5.59
$ws.Range("X4").Value = 4.15
$ws.Range("Y4").Value = 9.390000000000001
$ws.Range("Z4").Value = 6.84
$ws.Range("AA4").Value = 40.26
$ws.Range("AB4").Value = 2452.18
$ws.Range("AC4").Value = 1024
$ws.Range("AD4").Value = 19.42
$ws.Range("AE4").Value = 11431
$ws.Range("AF4").Value = 1.74
$ws.Range("AG4").Value = 350
$ws.Range("AH4").Value = 1.76
$ws.Range("AI4").Value = 33.57
$ws.Range("AJ4").Value = 17588764
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("V4").ClearContents()

# Row 5
$ws.Range("D5").Value = 3733
$ws.Range("E5").Value = 112
$ws.Range("F5").Value = 112
$ws.Range("G5").Value = 97
$ws.Range("H5").Value = 88
$ws.Range("I5").Value = 82
$ws.Range("J5").Value = 6
$ws.Range("K5").Value = 3722
$ws.Range("L5").Value = 676
$ws.Range("M5").Value = 3045
$ws.Range("N5").Value = 3034
$ws.Range("O5").Value = 11
$ws.Range("P5").Value = 124
$ws.Range("Q5").Value = 40
$ws.Range("R5").Value = -1103
$ws.Range("S5").Value = 1016
$ws.Range("T5").Value = 50
$ws.Range("U5").Value = -10
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 3.01
$ws.Range("X5").Value = 2.36
$ws.Range("Y5").Value = 3.27
$ws.Range("Z5").Value = 2.7
$ws.Range("AA5").Value = 22.21
$ws.Range("AB5").Value = 2372.93
$ws.Range("AC5").Value = 464
$ws.Range("AD5").Value = 35.81
$ws.Range("AE5").Value = 12354
$ws.Range("AF5").Value = 1.34
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 24823163

# Row 6
$ws.Range("D6").Value = 3455
$ws.Range("E6").Value = -190
$ws.Range("F6").Value = -190
$ws.Range("G6").Value = -143
$ws.Range("H6").Value = -117
$ws.Range("I6").Value = -119
$ws.Range("K6").Value = 3881
$ws.Range("L6").Value = 693
$ws.Range("M6").Value = 3188
$ws.Range("N6").Value = 3188
$ws.Range("P6").Value = 135
$ws.Range("Q6").Value = -395
$ws.Range("R6").Value = 306
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = 75
$ws.Range("U6").Value = -471
$ws.Range("V6").Value = 29
$ws.Range("W6").Value = -5.49
$ws.Range("X6").Value = -3.38
$ws.Range("Y6").Value = -3.83
$ws.Range("Z6").Value = -3.08
$ws.Range("AA6").Value = 21.73
$ws.Range("AB6").Value = 2277.83
$ws.Range("AC6").Value = -478
$ws.Range("AD6").Value = -22.78
$ws.Range("AE6").Value = 11909
$ws.Range("AF6").Value = 0.92
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 27033459
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()

# Row 7
$ws.Range("D7").Value = 3500
$ws.Range("E7").Value = -50
$ws.Range("W7").Value = -1.43
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").Value = 3370
$ws.Range("E8").Value = 80
$ws.Range("W8").Value = 2.37
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()

